# The commit adds one new daily price record for "Achicoria" (Primera
# quality, Región Metropolitana origin, 18-unit box) to the Vega Modelo
# de Temuco sheet. In the source data the rows are kept sorted by date,
# so the new record is inserted as row 52 and every following record
# shifts down by one row (old row 110 becomes the new row 111).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 52; Excel shifts rows 52:110 down
# to 53:111 automatically (dimension is updated by the engine as well).
$ws.Rows.Item(52).Insert()

$ws.Cells.Item(52, 1).Value  = 10
$ws.Cells.Item(52, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(52, 3).Value  = 'La Araucanía'
$ws.Cells.Item(52, 4).Value  = 45079
$ws.Cells.Item(52, 5).Value  = 9
$ws.Cells.Item(52, 6).Value  = 100112010
$ws.Cells.Item(52, 7).Value  = 'Achicoria'
$ws.Cells.Item(52, 8).Value  = 'Sin especificar'
$ws.Cells.Item(52, 9).Value  = 'Primera'
$ws.Cells.Item(52, 10).Value = 45
$ws.Cells.Item(52, 11).Value = 10000
$ws.Cells.Item(52, 12).Value = 10000
$ws.Cells.Item(52, 13).Value = 10000
$ws.Cells.Item(52, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(52, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(52, 16).Value = 556
$ws.Cells.Item(52, 17).Value = 18
$ws.Cells.Item(52, 18).Value = 'Hortaliza'

Write-Host "Inserted new row 52 (Achicoria, Región Metropolitana, 18 unidades)"
